$d = $word.ActiveDocument
$rng = $d.Content
$found = $rng.Find.Execute("<Student Name>                                                 < Date>", $false, $false, $false, $false, $false, $true, 1, $false, "", 1)
if (-not $found) {
  throw "Could not find target placeholder text"
}
$para = $rng.Paragraphs(1)
$fullRange = $d.Range($rng.Start, $para.Range.End)

$xml = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"><w:body><w:p w14:paraId="797DF8FE" w14:textId="77777777" w:rsidR="00362377" w:rsidRPr="00362377" w:rsidRDefault="00362377" w:rsidP="00362377"><w:pPr><w:autoSpaceDE w:val="0"/><w:autoSpaceDN w:val="0"/><w:adjustRightInd w:val="0"/><w:spacing w:after="200" w:line="276" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Verdana" w:eastAsia="Times New Roman" w:hAnsi="Verdana" w:cs="TimesNewRomanPS-BoldMT"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:eastAsia="Times New Roman" w:hAnsi="Verdana" w:cs="TimesNewRomanPS-BoldMT"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">Thomas </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:eastAsia="Times New Roman" w:hAnsi="Verdana" w:cs="TimesNewRomanPS-BoldMT"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t>Ziegelmann</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:rFonts w:ascii="Verdana" w:eastAsia="Times New Roman" w:hAnsi="Verdana" w:cs="TimesNewRomanPS-BoldMT"/><w:sz w:val="20"/><w:szCs w:val="20"/></w:rPr><w:t xml:space="preserve">                                                 January 1, 2025</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$fullRange.InsertXML($xml)
